$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 2.05
$ws.Range("Z2").Value = 9.5
$ws.Range("AC2").Value = 9
$ws.Range("AD2").Value = 8.5
$ws.Range("AI2").Value = 21
$ws.Range("AO2").Value = 7.5

$ws.Range("A3").Value = "QmuqFgzh"
$ws.Range("C3").Value = "20:30"
$ws.Range("E3").Value = "Racing Montevideo"
$ws.Range("F3").Value = "Nacional"
$ws.Range("G3").Value = 6.5
$ws.Range("H3").Value = 3.7
$ws.Range("I3").Value = 1.6
$ws.Range("J3").Value = 7
$ws.Range("K3").Value = 2.05
$ws.Range("L3").Value = 2.25
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.63
$ws.Range("Q3").Value = 2.4
$ws.Range("R3").Value = 1.53
$ws.Range("U3").Value = 2.5
$ws.Range("V3").Value = 1.5
$ws.Range("W3").Value = 12
$ws.Range("X3").Value = 29
$ws.Range("Y3").Value = 21
$ws.Range("Z3").Value = 81
$ws.Range("AA3").Value = 51
$ws.Range("AB3").Value = 67
$ws.Range("AD3").Value = 7.5
$ws.Range("AE3").Value = 26
$ws.Range("AF3").Value = 101
$ws.Range("AG3").Value = 5
$ws.Range("AH3").Value = 6
$ws.Range("AI3").Value = 9.5
$ws.Range("AJ3").Value = 11
$ws.Range("AK3").Value = 17
$ws.Range("AM3").Value = 351
$ws.Range("AN3").Value = 7.5
$ws.Range("AO3").Value = 41
$ws.Range("AP3").Value = 51
$ws.Range("AQ3").Value = 151
$ws.Range("AR3").Value = 251
$ws.Range("AS3").Value = 301
$ws.Range("AU3").Value = 10
$ws.Range("AV3").Value = 81
$ws.Range("AW3").Value = 3.25
$ws.Range("AX3").Value = 8.5
$ws.Range("AY3").Value = 26
$ws.Range("AZ3").Value = 29
$ws.Range("BA3").Value = 67
$ws.Range("BB3").Value = 251
